$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 29, shifting the existing rows 29-41 down to 30-42.
$ws.Rows(29).Insert()

# The Insert() call above copies formatting down from row 28 into the new
# row 29, which mints slightly different (border-less) style variants for
# E29/F29. Re-copy the real formatting from row 28 (C:F) onto row 29 so the
# new row matches the styles used by its sibling rows (no explicit style on
# C, style 21 on D, style 2 on E, style 3 on F).
$ws.Range("C28:F28").Copy()
$ws.Range("C29:F29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# New row 29: "total time v4" / "set cache = true" / 1.217
$ws.Range("C29").Value = "total time v4"
$ws.Range("D29").Value = "set cache = true"
$ws.Range("E29").Value = 1.2170000000000001
$ws.Rows(29).RowHeight = 17

# New row 30 (old row 29, shifted down): update the measured time.
$ws.Range("E30").Value = 1.766

# New row 32 (old blank row 31, shifted down): fill in the new data point
# "def arr vectors" / "cache = true" / 0.492, reusing row 31's C:D styling.
$ws.Range("C31:D31").Copy()
$ws.Range("C32:D32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("C32").Value = "def arr vectors"
$ws.Range("D32").Value = "cache = true"
$ws.Range("E32").Value = 0.49199999999999999
$ws.Rows(32).RowHeight = 17

# Update the tracked selection to match the author's last selection.
$ws.Range("I30").Select()
